$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-09 Wednesday", "2025-04-10 Thursday"),
    @("70×11=", "53×86="),
    @("33×45=", "45×79="),
    @("39×55=", "31×86="),
    @("83×21=", "60×30="),
    @("85×76=", "70×85="),
    @("75×18=", "21×65="),
    @("65×53=", "19×43="),
    @("86×78=", "71×20="),
    @("39×31=", "35×67="),
    @("57×69=", "27×91="),
    @("90×57=", "70×57="),
    @("65×35=", "83×16="),
    @("86×53=", "59×29="),
    @("82×77=", "58×36="),
    @("73×41=", "85×79="),
    @("37×29=", "93×22="),
    @("41×93=", "48×35="),
    @("66×95=", "36×96="),
    @("98×94=", "73×46="),
    @("62×23=", "25×57="),
    @("19×82=", "31×19="),
    @("93×76=", "46×49="),
    @("78×71=", "22×12="),
    @("55×54=", "89×55="),
    @("93×83=", "12×72=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
